# Highlight the "Merge Sort, Quick Sort" bullet (run text + paragraph
# mark) in cyan, matching the existing highlighting already applied to
# the preceding "Bubble Sort, Selection Sort, Insertion Sort" bullet.
#
# wdColorIndex 3 = wdTurquoise, which this engine serialises to
# <w:highlight w:val="cyan"/> in the OOXML.
$d = $word.ActiveDocument

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "Merge Sort, Quick Sort") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    # Going through .Font (rather than the Range directly) makes Word
    # apply the highlight to both the run(s) in the paragraph AND the
    # paragraph mark itself, matching how the sibling bullet above it
    # is already formatted.
    $target.Range.Font.HighlightColorIndex = 3
}
